$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Note: This schedule allows..." paragraph (right after the
#    "#4 - Platform Supervisor" / blank paragraph, before "End of phase 1").
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Note: This schedule allows for one extra week spare*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark that sits after "...embedded Linux".
# ---------------------------------------------------------------------------
try {
    $oldMark = $d.Bookmarks.Item("_GoBack")
    $oldMark.Delete()
} catch {
    # no existing _GoBack bookmark - nothing to remove
}

# ---------------------------------------------------------------------------
# 3) "End of phase 3 ... October 28th, 2016 - completion of EL" paragraph:
#      - "October 28" (bold) becomes "November 4" (not bold)
#      - append a new run " (Slip week used)" at the end of the paragraph
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("October 28") | Out-Null
$rng.Text = "November 4"
$rng.Font.Bold = 0

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "End of phase 3*completion of EL*") {
        $endRng = $d.Range($p.Range.End - 1, $p.Range.End - 1)
        $endRng.InsertAfter(" (Slip week used)")
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Remove the stray <w:lastRenderedPageBreak/> in front of
#    "Robot controlled over Wi-Fi through platform supervisor".
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("Robot controlled over Wi-Fi through platform supervisor", $true, $false, $false, $false, $false, $true, 1, $false, "Robot controlled over Wi-Fi through platform supervisor", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark to the very last (empty) paragraph of the
#    document body.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
